$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "M_PL" in column K, copying the header style (bold, centered, bordered)
# from the neighboring header cell J1 so it matches the existing header formatting.
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "M_PL"

# New profit (M_PL) values per region, keyed by row number.
$profits = @{
    2 = 106960237402
    3 = 137173931430
    4 = 42119558332
    5 = 9821205357
    6 = 885447038872
    7 = 12956669707
    8 = 3720464869
}

# For each data row, shift existing columns C..J right into D..K (read right-to-left so
# values aren't clobbered before they're copied), then place the new profit value in C.
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 11).Value = $ws.Cells.Item($r, 10).Value()
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 9).Value()
    $ws.Cells.Item($r, 9).Value = $ws.Cells.Item($r, 8).Value()
    $ws.Cells.Item($r, 8).Value = $ws.Cells.Item($r, 7).Value()
    $ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 6).Value()
    $ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 5).Value()
    $ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 4).Value()
    $ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 3).Value()
    $ws.Cells.Item($r, 3).Value = $profits[$r]
}
